$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new cell F1 value and style (no shifting of existing rows)
# Reuse the existing Times New Roman / size 12 font (fontId 2) by copying
# format from a cell that already uses it, then apply right alignment.
$src = $ws.Range("A13")
$cell = $ws.Range("F1")
$cell.Value = "Жадвал 22"
$src.Copy()
$cell.PasteSpecial(-4122)  # xlPasteFormats
$cell.HorizontalAlignment = -4152  # xlRight
$cell.VerticalAlignment = -4107    # xlGeneral (default, no explicit vertical)
$excel.CutCopyMode = $false

# Select A5:F5 range (the merged title cell) as the diff shows
$ws.Range("A5:F5").Select()

$wb.Save()
